# Update the "Last Updated: 20-Apr-20" date stamp to "21-Apr-20" on every
# slide, touching only the characters that actually changed so existing run
# formatting (bold headers, Calibri font, hyperlinks, line breaks, etc.) is
# left completely intact.

$p = $ppt.ActivePresentation

$oldText = "20-Apr-20"
$newText = "21-Apr-20"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)

        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text

            if ($full) {
                $idx = $full.IndexOf($oldText)
                while ($idx -ge 0) {
                    $sub = $tr.Characters($idx + 1, $oldText.Length)
                    $sub.Text = $newText

                    $full = $tr.Text
                    $idx = $full.IndexOf($oldText, $idx + $newText.Length)
                }
            }
        }
    }
}
